$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.982.09'
$ws.Range("E2").Value = '  +3.23%  '
$ws.Range("D3").Value = '3.193.97'
$ws.Range("E3").Value = '  +1.75%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.92'
$ws.Range("E6").Value = '  +4.09%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +2.86%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.34'
$ws.Range("E9").Value = '  +0.15%  '
$ws.Range("E10").Value = '  +4.26%  '
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("D12").Value = '3.745.60'
$ws.Range("E12").Value = '  +1.77%  '
$ws.Range("E13").Value = '  -1.49%  '
$ws.Range("E14").Value = '  +3.06%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.03'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").Value = '60.028.27'
$ws.Range("E16").Value = '  +3.09%  '
$ws.Range("D17").Value = '3.191.08'
$ws.Range("E17").Value = '  +1.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.22'
$ws.Range("E18").Value = '  -0.74%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.07'
$ws.Range("E19").Value = '  +0.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.35'
$ws.Range("E20").Value = '  +1.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.54'
$ws.Range("E21").Value = '  +1.50%  '
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").Value = '  +2.53%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.37'
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("E25").Value = '  +2.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.88'
$ws.Range("E26").Value = '  +11.36%  '
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("E28").Value = '  +2.20%  '
$ws.Range("E29").Value = '  +1.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.43'
$ws.Range("E30").Value = '  +3.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.19'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.23'
$ws.Range("E33").Value = '  +4.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.65'
$ws.Range("E34").Value = '  +5.34%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.08'
$ws.Range("E35").Value = '  -3.47%  '
$ws.Range("E36").Value = '  -0.20%  '
$ws.Range("D37").Value = '2.779.65'
$ws.Range("E37").Value = '  +5.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.74'
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0714'
$ws.Range("E39").Value = '  +5.57%  '
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.25'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.75'
$ws.Range("E42").Value = '  +1.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.729'
$ws.Range("E43").Value = '  +4.16%  '
$ws.Range("E44").Value = '  +4.99%  '
$ws.Range("D45").Value = '3.234.79'
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("E46").Value = '  +2.48%  '
$ws.Range("E47").Value = '  +1.22%  '
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.801'
$ws.Range("E49").Value = '  +6.72%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.58'
$ws.Range("E50").Value = '  +1.04%  '
$ws.Range("E51").Value = '  -0.03%  '
